$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that Word had left at the very start
#    of the document (first Heading1 paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Re-create the "_GoBack" bookmark further down, inside the title text
#    "ncview and ncBrowse" that introduces the Solution section (this is the
#    third occurrence of that phrase in the document - the first two are the
#    cover-page title and the "Aim:" sentence).
# ---------------------------------------------------------------------------
$needle = "ncview and ncBrowse"
$searchStart = 0
$lastMatch = $null
for ($i = 0; $i -lt 50; $i++) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $lastMatch = $d.Range($rng.Start, $rng.End)
    $searchStart = $rng.End
}

if ($lastMatch -ne $null) {
    # split point is right after "ncview and ncB" (14 characters in)
    $splitPos = $lastMatch.Start + 14
    $bmRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 3. Bold the "1. Let's look at the contents of an existing NetCDF file with
#    ncview." title line (including its paragraph mark).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("1. Let's look at the contents of an existing NetCDF file with")) {
        $p.Range.Font.Bold = 1
        break
    }
}

# ---------------------------------------------------------------------------
# 4. Bold the "2. Let's use ncBrowse to look at some agricultural emissions
#    data. " title line's visible text, while keeping the paragraph mark
#    itself explicitly not-bold.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("2. Let's use ncBrowse to look at some agricultural emissions data")) {
        $p.Range.Font.Bold = 0
        $textEnd = $p.Range.End - 1
        $textRange = $d.Range($p.Range.Start, $textEnd)
        $textRange.Font.Bold = 1
        break
    }
}

# ---------------------------------------------------------------------------
# 5. The built-in "Normal Table" style no longer needs to behave as a quick
#    style.
# ---------------------------------------------------------------------------
$tableStyle = $d.Styles.Item("Normal Table")
$tableStyle.QuickStyle = $false

Write-Host "edit complete"
